# Apply the OOXML changes described by the diff:
#  1. Every "datetimeFigureOut" date field (on the slide master and on
#     each of its 11 slide layouts) goes from 10/25/2021 -> 10/28/2021.
#  2. The "TextBox 7" shape ("Ligand annotations (k)") on slide 1 is
#     nudged to a new position (its size is unchanged).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update every date placeholder (ppPlaceholderDate = 16) found on the
#    slide master and all of its custom layouts.
# ---------------------------------------------------------------------
$oldDate = "10/25/2021"
$newDate = "10/28/2021"

$master = $p.SlideMaster
$targets = New-Object System.Collections.ArrayList
[void]$targets.Add($master)

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    [void]$targets.Add($layouts.Item($i))
}

foreach ($t in $targets) {
    for ($j = 1; $j -le $t.Shapes.Count; $j++) {
        $sh = $t.Shapes.Item($j)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {}
        if ($isDate) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Move "TextBox 7" on slide 1 to its new position.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
$box = $s.Shapes.Item(7)
if ($box.Name -eq "TextBox 7") {
    $box.Left = 6935760
    $box.Top = 6472717
}
